$wb = $excel.ActiveWorkbook

# --- Rename the "Include from RoleClass" sheet to "Include #0" ---
$wsInclude = $wb.Worksheets.Item("Include from RoleClass")
$wsInclude.Name = "Include #0"

# --- Update the Metadata sheet ---
$ws = $wb.Worksheets.Item("Metadata")

# Shift rows 11-14 (Description, Purpose, Copyright, Immutable) down one row
# to 12-15, making room for a new "Jurisdiction" property row at row 11.
for ($r = 14; $r -ge 11; $r--) {
    $dst = $r + 1
    $ws.Cells.Item($dst, 1).Value = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($dst, 2).Value = $ws.Cells.Item($r, 2).Value2
}

# The newly vacated last row (15) needs the same formatting as the row
# above it (14) since it previously held no data/format of its own.
$ws.Range("A14:B14").Copy()
$ws.Range("A15:B15").PasteSpecial(-4122)

# Populate the new "Jurisdiction" row.
$ws.Cells.Item(11, 1).Value = "Jurisdiction"
$ws.Cells.Item(11, 2).Value = ""

# Bump the Version and Date metadata values.
$ws.Cells.Item(3, 2).Value = "2.0.1-sd-202510-matchbox-patch"
$ws.Cells.Item(8, 2).Value = "2025-10-29T22:15:57+01:00"
